# Auto-generated Excel COM-interop script applying the cryptos.xlsx update
# Commit: Updated cryptos list on Mon May 22 22:47:17 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.011.80"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.02"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.56"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4631"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3703"
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07341"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8768"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07895"
$ws.Range("E11").Value = "  +3.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.70"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.859.97"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.337"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.52"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008835"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.79"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.040.18"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.103"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.52"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.050.66"
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.037"
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.66"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08866"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.961"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7308"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.439"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.471"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.076"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05217"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.958"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.103"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5158"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1627"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.151"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4828"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.19"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.93"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.626"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06206"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.77"
$ws.Range("E51").Value = "  +0.17%  "
